# Updates the cryptos worksheet with refreshed price/volume snapshot data.
# D-column price strings contain literal "." grouping characters and
# would otherwise be re-interpreted by Excel as numbers; the leading
# apostrophe forces them to remain text, matching the source data feed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.176.50"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "'1.902.89"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'306.33"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "'0.5253"
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("D8").Value = "'0.3775"
$ws.Range("E8").Value = "  +1.30%  "
$ws.Range("D9").Value = "'0.07264"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("E10").Value = "  +0.12%  "
$ws.Range("D11").Value = "'0.8989"
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "'0.08383"
$ws.Range("E12").Value = "  +9.50%  "
$ws.Range("D13").Value = "'1.880.09"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'94.84"
$ws.Range("E14").Value = "  -0.35%  "
$ws.Range("D15").Value = "'5.272"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").Value = "'0.000008610"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'27.207.48"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "'5.062"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "'2.134.26"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Value = "'6.436"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "'146.76"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Value = "'2.279"
$ws.Range("E26").Value = "  +5.83%  "
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").Value = "'18.15"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "'114.91"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("D30").Value = "'4.930"
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'4.792"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "'0.09286"
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Value = "'0.8156"
$ws.Range("E33").Value = "  +6.87%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("D35").Value = "'1.238"
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D36").Value = "'2.951"
$ws.Range("E36").Value = "  -2.39%  "
$ws.Range("D37").Value = "'3.349"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "'2.588"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").Value = "'0.5715"
$ws.Range("E39").Value = "  +1.57%  "
$ws.Range("E40").Value = "  -0.56%  "
$ws.Range("D41").Value = "'1.070"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").Value = "'8.954"
$ws.Range("E43").Value = "  +0.98%  "
$ws.Range("D44").Value = "'117.97"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("D45").Value = "'0.1514"
$ws.Range("E45").Value = "  +0.39%  "
$ws.Range("D46").Value = "'0.4835"
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("D47").Value = "'10.22"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  +0.24%  "
$ws.Range("D49").Value = "'1.618"
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "'37.49"
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("D51").Value = "'63.64"
$ws.Range("E51").Value = "  +0.09%  "
